# "Generate Report for handoff"
# - Two new source files (24c887d2-... and 978a7524-...) show up as "Ready for handoff"
#   on all three sheets, inserted right before the ".localization-config" row.
# - The two pre-existing in-flight files (6802e21e-... and dda88031-...) move from
#   "Ready for handoff" to "In Translation" on the Overview sheet status columns.
# - zh-cn / de-de detail sheets get matching new rows with handoff file + datetime.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

# Make room: push the ".localization-config" row from row 4 down to row 6.
$ws.Rows("4:5").Insert()

# Duplicate row 3's formatting (hyperlink-style A column, plain B/C) onto the
# two freshly inserted rows so fonts/number formats line up with the rest of
# the table before we overwrite the cell contents.
$ws.Rows("3:3").Copy()
$ws.Paste($ws.Rows("4:4"))
$ws.Rows("3:3").Copy()
$ws.Paste($ws.Rows("5:5"))
$excel.CutCopyMode = $false

# Existing in-flight files: status flips to "In Translation".
$ws.Range("B2").Value2 = "In Translation"
$ws.Range("C2").Value2 = "In Translation"
$ws.Range("B3").Value2 = "In Translation"
$ws.Range("C3").Value2 = "In Translation"

# New row 4: 24c887d2-...
$ws.Range("A4").Value2 = "24c887d2-f6c7-4377-b1a1-20baaee63615.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "Ready for handoff"

# New row 5: 978a7524-...
$ws.Range("A5").Value2 = "978a7524-7d3a-4a07-9038-42538f5b36e1.md"
$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("C5").Value2 = "Ready for handoff"

# Row 6 (the old row 4) keeps its original text; only its hyperlink target
# needs to be re-pointed below since it physically moved rows.

# The inserted rows left the hyperlink collection stale (the old row 4 link
# now points at a blank cell), so rebuild every link on this sheet in display
# order and restore the hyperlink look (underline, cornflower-blue) afterward.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/6802e21e-fecc-4054-a3e0-421fa9b0cead.md", "", "", "6802e21e-fecc-4054-a3e0-421fa9b0cead.md")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/dda88031-d710-4511-b347-eacbd39859af.md", "", "", "dda88031-d710-4511-b347-eacbd39859af.md")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/24c887d2-f6c7-4377-b1a1-20baaee63615.md", "", "", "24c887d2-f6c7-4377-b1a1-20baaee63615.md")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/978a7524-7d3a-4a07-9038-42538f5b36e1.md", "", "", "978a7524-7d3a-4a07-9038-42538f5b36e1.md")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/.localization-config", "", "", ".localization-config")

$r = $ws.Range("A2:A6")
$r.Font.Underline = 2
$r.Font.Color = 15570276

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn detail
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Rows("4:5").Insert()
$ws.Rows("3:3").Copy()
$ws.Paste($ws.Rows("4:4"))
$ws.Rows("3:3").Copy()
$ws.Paste($ws.Rows("5:5"))
$excel.CutCopyMode = $false

$ws.Range("B2").Value2 = "In Translation"
$ws.Range("B3").Value2 = "In Translation"

$ws.Range("A4").Value2 = "24c887d2-f6c7-4377-b1a1-20baaee63615.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "24c887d2-f6c7-4377-b1a1-20baaee63615.bd4d7f3deea8eb78deaea883ac5a9a549f191bc5.zh-cn.xlf"
$ws.Range("D4").Value2 = "2016-01-25 10:36:59"
$ws.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws.Range("H4").Value2 = "Include"

$ws.Range("A5").Value2 = "978a7524-7d3a-4a07-9038-42538f5b36e1.md"
$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("C5").Value2 = "978a7524-7d3a-4a07-9038-42538f5b36e1.e1fea881e098dc711eac17e4e7d6bd425e6dda60.zh-cn.xlf"
$ws.Range("D5").Value2 = "2016-01-25 10:36:59"
$ws.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws.Range("H5").Value2 = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/6802e21e-fecc-4054-a3e0-421fa9b0cead.md", "", "", "6802e21e-fecc-4054-a3e0-421fa9b0cead.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9946747fbd241a866e38a45a37873f1e875bb1d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/6802e21e-fecc-4054-a3e0-421fa9b0cead.0abc74fb7629d3542d0e92f888702fcfa0a680bd.zh-cn.xlf", "", "", "6802e21e-fecc-4054-a3e0-421fa9b0cead.0abc74fb7629d3542d0e92f888702fcfa0a680bd.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/dda88031-d710-4511-b347-eacbd39859af.md", "", "", "dda88031-d710-4511-b347-eacbd39859af.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9946747fbd241a866e38a45a37873f1e875bb1d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/dda88031-d710-4511-b347-eacbd39859af.e29ce5550d52216bf06ae615210c8417b4e22284.zh-cn.xlf", "", "", "dda88031-d710-4511-b347-eacbd39859af.e29ce5550d52216bf06ae615210c8417b4e22284.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/24c887d2-f6c7-4377-b1a1-20baaee63615.md", "", "", "24c887d2-f6c7-4377-b1a1-20baaee63615.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9946747fbd241a866e38a45a37873f1e875bb1d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/24c887d2-f6c7-4377-b1a1-20baaee63615.bd4d7f3deea8eb78deaea883ac5a9a549f191bc5.zh-cn.xlf", "", "", "24c887d2-f6c7-4377-b1a1-20baaee63615.bd4d7f3deea8eb78deaea883ac5a9a549f191bc5.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/978a7524-7d3a-4a07-9038-42538f5b36e1.md", "", "", "978a7524-7d3a-4a07-9038-42538f5b36e1.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9946747fbd241a866e38a45a37873f1e875bb1d7/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/978a7524-7d3a-4a07-9038-42538f5b36e1.e1fea881e098dc711eac17e4e7d6bd425e6dda60.zh-cn.xlf", "", "", "978a7524-7d3a-4a07-9038-42538f5b36e1.e1fea881e098dc711eac17e4e7d6bd425e6dda60.zh-cn.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/.localization-config", "", "", ".localization-config")

$r = $ws.Range("A2:A6")
$r.Font.Underline = 2
$r.Font.Color = 15570276
$r = $ws.Range("C2:C5")
$r.Font.Underline = 2
$r.Font.Color = 15570276

# ---------------------------------------------------------------------------
# Sheet 3: de-de detail
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Rows("4:5").Insert()
$ws.Rows("3:3").Copy()
$ws.Paste($ws.Rows("4:4"))
$ws.Rows("3:3").Copy()
$ws.Paste($ws.Rows("5:5"))
$excel.CutCopyMode = $false

$ws.Range("B2").Value2 = "In Translation"
$ws.Range("B3").Value2 = "In Translation"

$ws.Range("A4").Value2 = "24c887d2-f6c7-4377-b1a1-20baaee63615.md"
$ws.Range("B4").Value2 = "Ready for handoff"
$ws.Range("C4").Value2 = "24c887d2-f6c7-4377-b1a1-20baaee63615.bd4d7f3deea8eb78deaea883ac5a9a549f191bc5.de-de.xlf"
$ws.Range("D4").Value2 = "2016-01-25 10:37:10"
$ws.Range("G4").Value2 = "0001-01-01 00:00:00"
$ws.Range("H4").Value2 = "Include"

$ws.Range("A5").Value2 = "978a7524-7d3a-4a07-9038-42538f5b36e1.md"
$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("C5").Value2 = "978a7524-7d3a-4a07-9038-42538f5b36e1.e1fea881e098dc711eac17e4e7d6bd425e6dda60.de-de.xlf"
$ws.Range("D5").Value2 = "2016-01-25 10:37:10"
$ws.Range("G5").Value2 = "0001-01-01 00:00:00"
$ws.Range("H5").Value2 = "Include"

$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/6802e21e-fecc-4054-a3e0-421fa9b0cead.md", "", "", "6802e21e-fecc-4054-a3e0-421fa9b0cead.md")
$ws.Hyperlinks.Add($ws.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/baa3fa05ffe18f9193b6a08928e68ecedc9feac6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/6802e21e-fecc-4054-a3e0-421fa9b0cead.0abc74fb7629d3542d0e92f888702fcfa0a680bd.de-de.xlf", "", "", "6802e21e-fecc-4054-a3e0-421fa9b0cead.0abc74fb7629d3542d0e92f888702fcfa0a680bd.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/dda88031-d710-4511-b347-eacbd39859af.md", "", "", "dda88031-d710-4511-b347-eacbd39859af.md")
$ws.Hyperlinks.Add($ws.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/baa3fa05ffe18f9193b6a08928e68ecedc9feac6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/dda88031-d710-4511-b347-eacbd39859af.e29ce5550d52216bf06ae615210c8417b4e22284.de-de.xlf", "", "", "dda88031-d710-4511-b347-eacbd39859af.e29ce5550d52216bf06ae615210c8417b4e22284.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/24c887d2-f6c7-4377-b1a1-20baaee63615.md", "", "", "24c887d2-f6c7-4377-b1a1-20baaee63615.md")
$ws.Hyperlinks.Add($ws.Range("C4"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/baa3fa05ffe18f9193b6a08928e68ecedc9feac6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/24c887d2-f6c7-4377-b1a1-20baaee63615.bd4d7f3deea8eb78deaea883ac5a9a549f191bc5.de-de.xlf", "", "", "24c887d2-f6c7-4377-b1a1-20baaee63615.bd4d7f3deea8eb78deaea883ac5a9a549f191bc5.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A5"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/e2e/978a7524-7d3a-4a07-9038-42538f5b36e1.md", "", "", "978a7524-7d3a-4a07-9038-42538f5b36e1.md")
$ws.Hyperlinks.Add($ws.Range("C5"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/baa3fa05ffe18f9193b6a08928e68ecedc9feac6/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/978a7524-7d3a-4a07-9038-42538f5b36e1.e1fea881e098dc711eac17e4e7d6bd425e6dda60.de-de.xlf", "", "", "978a7524-7d3a-4a07-9038-42538f5b36e1.e1fea881e098dc711eac17e4e7d6bd425e6dda60.de-de.xlf")
$ws.Hyperlinks.Add($ws.Range("A6"), "https://github.com/OpenLocalizationTest/oltest/blob/f173daa544d8366040bf14e77ea7cc48ee9722d8/.localization-config", "", "", ".localization-config")

$r = $ws.Range("A2:A6")
$r.Font.Underline = 2
$r.Font.Color = 15570276
$r = $ws.Range("C2:C5")
$r.Font.Underline = 2
$r.Font.Color = 15570276

$ws = $wb.Worksheets.Item("Overview")
$ws.Activate()

Write-Host "edit complete"
